$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 334000000
$ws.Range("J70").Value = 334000000
$ws.Range("L70").Value = 1002000000
$ws.Range("N70").Value = -1002000540

$ws.Range("H73").Value = 334000000
$ws.Range("J73").Value = 334000000
$ws.Range("L73").Value = 1002000000
$ws.Range("N73").Value = -1002001872

$ws.Range("H112").Value = 3321.2363
$ws.Range("I112").Value = 1332
$ws.Range("J112").Value = 3436
$ws.Range("K112").Value = 3996
$ws.Range("L112").Value = 10308
$ws.Range("M112").Value = -2888
$ws.Range("N112").Value = -12524

$ws.Range("H141").Value = 1333
$ws.Range("I141").Value = 1333
$ws.Range("K141").Value = 3999
$ws.Range("M141").Value = 1181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8424.549000000001
$ws.Range("I32").Value = 1882
$ws.Range("K32").Value = 1882
$ws.Range("M32").Value = -1595

$ws.Range("H61").Value = 6657.8887
$ws.Range("I61").Value = 4990.25
$ws.Range("K61").Value = 4990.25
$ws.Range("M61").Value = -4778.25

$ws.Range("H63").Value = 7399.2666
$ws.Range("J63").Value = 9299.299999999999
$ws.Range("L63").Value = 9299.299999999999
$ws.Range("N63").Value = -10671.3

$ws.Range("H66").Value = 7399.2666
$ws.Range("J66").Value = 9299.299999999999
$ws.Range("L66").Value = 46496.5
$ws.Range("N66").Value = -53360.5

$ws.Range("H136").Value = 6657.8887
$ws.Range("I136").Value = 4990.25
$ws.Range("K136").Value = 14970.75
$ws.Range("M136").Value = -12420.75

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 844
$ws.Range("I80").Value = 141
$ws.Range("K80").Value = 141
$ws.Range("M80").Value = 857

$ws.Range("H83").Value = 844
$ws.Range("I83").Value = 141
$ws.Range("K83").Value = 705
$ws.Range("M83").Value = 4287

$ws.Range("H94").Value = 15641723
$ws.Range("I94").Value = 31251146
$ws.Range("J94").Value = 32300.875
$ws.Range("K94").Value = 31251146
$ws.Range("L94").Value = 32300.875
$ws.Range("M94").Value = -31250695
$ws.Range("N94").Value = -33202.875

$ws.Range("H99").Value = 3036.35
$ws.Range("I99").Value = 3300.9412
$ws.Range("K99").Value = 3300.9412
$ws.Range("M99").Value = -1802.9412

$ws.Range("H107").Value = 3482.6
$ws.Range("I107").Value = 3042.0625
$ws.Range("J107").Value = 5244.75
$ws.Range("K107").Value = 3042.0625
$ws.Range("L107").Value = 5244.75
$ws.Range("M107").Value = -1122.0625
$ws.Range("N107").Value = -9084.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5718.959
$ws.Range("I31").Value = 4391.516
$ws.Range("J31").Value = 8005.1113
$ws.Range("K31").Value = 4391.516
$ws.Range("L31").Value = 8005.1113
$ws.Range("M31").Value = -4096.516
$ws.Range("N31").Value = -8595.1113

$ws.Range("H34").Value = 5718.959
$ws.Range("I34").Value = 4391.516
$ws.Range("J34").Value = 8005.1113
$ws.Range("K34").Value = 4391.516
$ws.Range("L34").Value = 8005.1113
$ws.Range("M34").Value = -4189.516
$ws.Range("N34").Value = -8409.1113

$ws.Range("H58").Value = 6818.75
$ws.Range("I58").Value = 7499.3335
$ws.Range("K58").Value = 7499.3335
$ws.Range("M58").Value = -7296.3335

$ws.Range("H99").Value = 5134.125
$ws.Range("I99").Value = 4391.643
$ws.Range("J99").Value = 6173.6
$ws.Range("K99").Value = 4391.643
$ws.Range("L99").Value = 6173.6
$ws.Range("M99").Value = -2893.643
$ws.Range("N99").Value = -9169.6

$ws.Range("H107").Value = 25001028
$ws.Range("I107").Value = 33334080
$ws.Range("K107").Value = 33334080
$ws.Range("M107").Value = -33332160

$ws.Range("H126").Value = 5134.125
$ws.Range("I126").Value = 4391.643
$ws.Range("J126").Value = 6173.6
$ws.Range("K126").Value = 13174.929
$ws.Range("L126").Value = 18520.8
$ws.Range("M126").Value = -10704.929
$ws.Range("N126").Value = -23460.8

$ws.Range("H136").Value = 6818.75
$ws.Range("I136").Value = 7499.3335
$ws.Range("K136").Value = 22498.0005
$ws.Range("M136").Value = -19948.0005

$ws.Range("H140").Value = 64591.8
$ws.Range("J140").Value = 64591.8
$ws.Range("L140").Value = 64591.8
$ws.Range("N140").Value = -74951.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 9998
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 9998
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 29994
$ws.Range("N32").Value = -30560
$ws.Range("M32").ClearContents()

$ws.Range("H33").Value = 96.44444
$ws.Range("I33").Value = 33.6
$ws.Range("J33").Value = 175
$ws.Range("K33").Value = 201.6
$ws.Range("L33").Value = 1050
$ws.Range("M33").Value = 81.39999999999998
$ws.Range("N33").Value = -1616

$ws.Range("H41").Value = 4950
$ws.Range("I41").Value = 4900
$ws.Range("K41").Value = 14700
$ws.Range("M41").Value = -14362

$ws.Range("H44").Value = 2142.5715
$ws.Range("I44").Value = 499.5
$ws.Range("K44").Value = 1498.5
$ws.Range("M44").Value = -1100.5

$ws.Range("H69").Value = 1033.3334
$ws.Range("I69").Value = 1033.3334
$ws.Range("K69").Value = 3100.0002
$ws.Range("M69").Value = -2289.0002

$ws.Range("H72").Value = 1033.3334
$ws.Range("I72").Value = 1033.3334
$ws.Range("K72").Value = 9300.000599999999
$ws.Range("M72").Value = -5244.000599999999

$ws.Range("H98").Value = 687.7143
$ws.Range("J98").Value = 687.7143
$ws.Range("L98").Value = 2063.1429
$ws.Range("N98").Value = -5059.1429

$ws.Range("H107").Value = 591.7895
$ws.Range("I107").Value = 552.4286
$ws.Range("J107").Value = 614.75
$ws.Range("K107").Value = 1657.2858
$ws.Range("L107").Value = 1844.25
$ws.Range("M107").Value = 262.7142000000001
$ws.Range("N107").Value = -5684.25

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 541.5
$ws.Range("I132").Value = 541.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4873.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2343.5
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 4804.467
$ws.Range("I134").Value = 1022.25
$ws.Range("K134").Value = 3066.75
$ws.Range("M134").Value = 2003.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5543.067
$ws.Range("I80").Value = 3843.5
$ws.Range("J80").Value = 7485.4287
$ws.Range("K80").Value = 3843.5
$ws.Range("L80").Value = 7485.4287
$ws.Range("M80").Value = -2845.5
$ws.Range("N80").Value = -9481.4287

$ws.Range("H83").Value = 5543.067
$ws.Range("I83").Value = 3843.5
$ws.Range("J83").Value = 7485.4287
$ws.Range("K83").Value = 19217.5
$ws.Range("L83").Value = 37427.14350000001
$ws.Range("M83").Value = -14225.5
$ws.Range("N83").Value = -47411.14350000001

$ws.Range("H122").Value = 859.8
$ws.Range("I122").Value = 859.8
$ws.Range("K122").Value = 2579.4
$ws.Range("M122").Value = -129.3999999999996

$ws.Range("H126").Value = 7523.36
$ws.Range("I126").Value = 6803.65
$ws.Range("K126").Value = 20410.95
$ws.Range("M126").Value = -17940.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4262
$ws.Range("I7").Value = 4294.5454
$ws.Range("K7").Value = 4294.5454
$ws.Range("M7").Value = -4182.5454

$ws.Range("H16").Value = 1064.762
$ws.Range("I16").Value = 695.6429000000001
$ws.Range("J16").Value = 1803
$ws.Range("K16").Value = 695.6429000000001
$ws.Range("L16").Value = 1803
$ws.Range("M16").Value = -525.6429000000001
$ws.Range("N16").Value = -2143

$ws.Range("H22").Value = 2256.8572
$ws.Range("I22").Value = 1825
$ws.Range("K22").Value = 1825
$ws.Range("M22").Value = -1530

$ws.Range("H27").Value = 2256.8572
$ws.Range("I27").Value = 1825
$ws.Range("K27").Value = 1825
$ws.Range("M27").Value = -1718

$ws.Range("H40").Value = 8690.362999999999
$ws.Range("I40").Value = 8690.362999999999
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 8690.362999999999
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -8554.362999999999
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 1967
$ws.Range("I46").Value = 1578.5
$ws.Range("J46").Value = 2277.8
$ws.Range("K46").Value = 1578.5
$ws.Range("L46").Value = 2277.8
$ws.Range("M46").Value = -1390.5
$ws.Range("N46").Value = -2653.8

$ws.Range("H55").Value = 119.4
$ws.Range("I55").Value = 100
$ws.Range("K55").Value = 100
$ws.Range("M55").Value = 73

$ws.Range("H126").Value = 4262
$ws.Range("I126").Value = 4294.5454
$ws.Range("K126").Value = 12883.6362
$ws.Range("M126").Value = -10413.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H62").Value = 9221.5
$ws.Range("I62").Value = 8962
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 8962
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -8338
$ws.Range("N62").Value = -11248

$ws.Range("H65").Value = 9221.5
$ws.Range("I65").Value = 8962
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 44810
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -41690
$ws.Range("N65").Value = -56240
